$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.195.34"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3
$ws.Range("D3").Value = "1.855.33"
$ws.Range("E3").Value = "  -0.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07766"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3068"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.49%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.21%  "

# Row 12
$ws.Range("D12").Value = "1.863.03"
$ws.Range("E12").Value = "  +0.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "91.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6857"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.39%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.522"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.44%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008459"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.20%  "

# Row 18
$ws.Range("D18").Value = "29.200.15"
$ws.Range("E18").Value = "  +0.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "

# Row 20
$ws.Range("D20").Value = "2.104.99"
$ws.Range("E20").Value = "  -0.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.527"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.36%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1500"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.34%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.851"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.37%  "

# Row 28
$ws.Range("E28").Value = "  -1.81%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.554"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.248"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.14%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.204"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.27%  "

# Row 32
$ws.Range("E32").Value = "  -1.23%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05228"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7597"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.44%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.169"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.38%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.840"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.75%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.708"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01864"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "

# Row 39
$ws.Range("D39").Value = "1.224.46"
$ws.Range("E39").Value = "  -1.30%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.731"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.69%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9992"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.004.28"
$ws.Range("E44").Value = "  -0.10%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.499"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.66%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.19%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5181"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.527"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.750"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.59%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.044"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
